$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.05172844094803208
$ws.Range("C2").Value = 0.9984787548029618
$ws.Range("D2").Value = 0.1671337180095297
$ws.Range("G2").Value = 0.3552643978832445
$ws.Range("H2").Value = 0.9990000000000001

# Row 3
$ws.Range("B3").Value = 0.06091766037437762
$ws.Range("C3").Value = 0.9994221686912212
$ws.Range("D3").Value = 0.1791403288070286
$ws.Range("G3").Value = 0.3552643978832445
$ws.Range("H3").Value = 0.9990000000000001

# Row 4
$ws.Range("B4").Value = 0.05602327423319848
$ws.Range("C4").Value = 0.9992414794489526
$ws.Range("D4").Value = 0.1895399969410011
$ws.Range("G4").Value = 0.3552643978832445
$ws.Range("H4").Value = 0.9990000000000001
